$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 671.4138
$ws.Range("I28").Value = 457.55
$ws.Range("J28").Value = 1146.6666
$ws.Range("K28").Value = 457.55
$ws.Range("L28").Value = 1146.6666
$ws.Range("M28").Value = 27.44999999999999
$ws.Range("N28").Value = -2116.6666
$ws.Range("H62").Value = 16812.625
$ws.Range("I62").Value = 30813.715
$ws.Range("J62").Value = 5922.8887
$ws.Range("K62").Value = 30813.715
$ws.Range("L62").Value = 5922.8887
$ws.Range("M62").Value = -30189.715
$ws.Range("N62").Value = -7170.8887
$ws.Range("H65").Value = 16812.625
$ws.Range("I65").Value = 30813.715
$ws.Range("J65").Value = 5922.8887
$ws.Range("K65").Value = 154068.575
$ws.Range("L65").Value = 29614.4435
$ws.Range("M65").Value = -150948.575
$ws.Range("N65").Value = -35854.4435
$ws.Range("H86").Value = 4614.5
$ws.Range("I86").Value = 2687.6667
$ws.Range("J86").Value = 6837.769
$ws.Range("K86").Value = 2687.6667
$ws.Range("L86").Value = 6837.769
$ws.Range("M86").Value = -1564.6667
$ws.Range("N86").Value = -9083.769
$ws.Range("H89").Value = 4614.5
$ws.Range("I89").Value = 2687.6667
$ws.Range("J89").Value = 6837.769
$ws.Range("K89").Value = 13438.3335
$ws.Range("L89").Value = 34188.845
$ws.Range("M89").Value = -7822.333500000001
$ws.Range("N89").Value = -45420.845
$ws.Range("H113").Value = 3708
$ws.Range("I113").Value = 3115.111
$ws.Range("J113").Value = 3974.8
$ws.Range("K113").Value = 3115.111
$ws.Range("L113").Value = 3974.8
$ws.Range("M113").Value = 138.8890000000001
$ws.Range("N113").Value = -10482.8
$ws.Range("H132").Value = 3415.6316
$ws.Range("I132").Value = 2166.389
$ws.Range("J132").Value = 4539.95
$ws.Range("K132").Value = 6499.167
$ws.Range("L132").Value = 13619.85
$ws.Range("M132").Value = -3969.167
$ws.Range("N132").Value = -18679.85
$ws.Range("H138").Value = 1595.8657
$ws.Range("I138").Value = 691.3929000000001
$ws.Range("J138").Value = 2245.2307
$ws.Range("K138").Value = 2074.1787
$ws.Range("L138").Value = 6735.6921
$ws.Range("M138").Value = 3065.8213
$ws.Range("N138").Value = -17015.6921
$ws.Range("H141").Value = 2010.4286
$ws.Range("I141").Value = 1512.1666
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 4536.4998
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = 643.5002000000004
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2326.6428
$ws.Range("I2").Value = 2463.3333
$ws.Range("K2").Value = 2463.3333
$ws.Range("M2").Value = -2350.3333
$ws.Range("H45").Value = 1234.16
$ws.Range("I45").Value = 954.625
$ws.Range("J45").Value = 1731.1111
$ws.Range("K45").Value = 954.625
$ws.Range("L45").Value = 1731.1111
$ws.Range("M45").Value = -577.625
$ws.Range("N45").Value = -2485.1111
$ws.Range("H61").Value = 330650.25
$ws.Range("I61").Value = 304995.97
$ws.Range("J61").Value = 360885.66
$ws.Range("K61").Value = 304995.97
$ws.Range("L61").Value = 360885.66
$ws.Range("M61").Value = -304783.97
$ws.Range("N61").Value = -361309.66
$ws.Range("H110").Value = 2340.9092
$ws.Range("I110").Value = 2431.25
$ws.Range("J110").Value = 2100
$ws.Range("K110").Value = 2431.25
$ws.Range("L110").Value = 2100
$ws.Range("M110").Value = -386.25
$ws.Range("N110").Value = -6190
$ws.Range("H116").Value = 2326.6428
$ws.Range("I116").Value = 2463.3333
$ws.Range("K116").Value = 2463.3333
$ws.Range("M116").Value = -169.3332999999998
$ws.Range("H136").Value = 330650.25
$ws.Range("I136").Value = 304995.97
$ws.Range("J136").Value = 360885.66
$ws.Range("K136").Value = 914987.9099999999
$ws.Range("L136").Value = 1082656.98
$ws.Range("M136").Value = -912437.9099999999
$ws.Range("N136").Value = -1087756.98

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2326.6428
$ws.Range("I3").Value = 2463.3333
$ws.Range("K3").Value = 2463.3333
$ws.Range("M3").Value = -2349.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5053693
$ws.Range("I62").Value = 7938891
$ws.Range("J62").Value = 4597
$ws.Range("K62").Value = 7938891
$ws.Range("L62").Value = 4597
$ws.Range("M62").Value = -7938267
$ws.Range("N62").Value = -5845
$ws.Range("H65").Value = 5053693
$ws.Range("I65").Value = 7938891
$ws.Range("J65").Value = 4597
$ws.Range("K65").Value = 39694455
$ws.Range("L65").Value = 22985
$ws.Range("M65").Value = -39691335
$ws.Range("N65").Value = -29225
$ws.Range("H134").Value = 2245.9666
$ws.Range("I134").Value = 1592.7858
$ws.Range("J134").Value = 2817.5
$ws.Range("K134").Value = 4778.357400000001
$ws.Range("L134").Value = 8452.5
$ws.Range("M134").Value = -2243.357400000001
$ws.Range("N134").Value = -13522.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 2783.9333
$ws.Range("I117").Value = 5981.8
$ws.Range("J117").Value = 1185
$ws.Range("K117").Value = 17945.4
$ws.Range("L117").Value = 3555
$ws.Range("M117").Value = -14503.4
$ws.Range("N117").Value = -10439
$ws.Range("H121").Value = 2292.7568
$ws.Range("I121").Value = 3285
$ws.Range("J121").Value = 1925.2593
$ws.Range("K121").Value = 9855
$ws.Range("L121").Value = 5775.7779
$ws.Range("M121").Value = -8545
$ws.Range("N121").Value = -8395.777900000001
$ws.Range("H140").Value = 2514.1667
$ws.Range("I140").Value = 2435.5557
$ws.Range("J140").Value = 2750
$ws.Range("K140").Value = 7306.6671
$ws.Range("L140").Value = 8250
$ws.Range("M140").Value = -2126.6671
$ws.Range("N140").Value = -18610

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 4626.5
$ws.Range("I21").Value = 4835.3335
$ws.Range("K21").Value = 4835.3335
$ws.Range("M21").Value = -4662.3335
$ws.Range("H30").Value = 4626.5
$ws.Range("I30").Value = 4835.3335
$ws.Range("K30").Value = 4835.3335
$ws.Range("M30").Value = -4730.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4281.784
$ws.Range("I136").Value = 2600.724
$ws.Range("J136").Value = 6497.727
$ws.Range("K136").Value = 7802.172
$ws.Range("L136").Value = 19493.181
$ws.Range("M136").Value = -5252.172
$ws.Range("N136").Value = -24593.181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 484.10526
$ws.Range("I107").Value = 450.2857
$ws.Range("J107").Value = 503.83334
$ws.Range("K107").Value = 1350.8571
$ws.Range("L107").Value = 1511.50002
$ws.Range("M107").Value = 569.1428999999998
$ws.Range("N107").Value = -5351.500019999999
$ws.Range("H113").Value = 307.2963
$ws.Range("I113").Value = 291.88
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 875.64
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1294.36
$ws.Range("N113").Value = -5840
$ws.Range("H132").Value = 1749.2616
$ws.Range("I132").Value = 1173.8518
$ws.Range("J132").Value = 2158.1052
$ws.Range("K132").Value = 3521.5554
$ws.Range("L132").Value = 6474.3156
$ws.Range("M132").Value = -991.5553999999997
$ws.Range("N132").Value = -11534.3156

Write-Host "Applied all changes"